$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Right-hand "binary" table (columns N:W), first block: rows 14-21 ---
# Column W used to CONCAT(T,",") (T21 stood alone as "=T21"); change it to
# mirror column K's pattern: CONCATENATE("B", N, O, P, Q, R, ",")
$ws.Range("W14").Formula = '=CONCATENATE("B",N14,O14,P14,Q14,R14,",")'
for ($r = 15; $r -le 21; $r++) {
    $ws.Cells.Item($r, 23).Formula = "=CONCATENATE(""B"",N$r,O$r,P$r,Q$r,R$r,"","")"
}

# --- Right-hand "binary" table, second block: rows 26-33 ---
$ws.Range("W26").Formula = '=CONCATENATE("B",N26,O26,P26,Q26,R26,",")'
for ($r = 27; $r -le 33; $r++) {
    $ws.Cells.Item($r, 23).Formula = "=CONCATENATE(""B"",N$r,O$r,P$r,Q$r,R$r,"","")"
}

# --- Move the active selection to T34 ---
$ws.Range("T34").Select() | Out-Null
